$wb = $excel.ActiveWorkbook

# ---- Sheet: Overview ----
$ws = $wb.Worksheets.Item('Overview')

# Update cell values (row 2 <-> row 3 content swap + refreshed status/dates)
$ws.Range('A2').Value = '77232830-5d71-4781-bf09-c9d381f516af.md'
$ws.Range('B2').Value = 'Handed back: in sync with en-US'
$ws.Range('C2').Value = 'Handed back: in sync with en-US'
$ws.Range('D2').Value = '2016-03-25 10:54:50'
$ws.Range('A3').Value = '2e932acd-e47f-4f3c-8372-e61745a5bd03.md'
$ws.Range('B3').Value = 'Ready for handoff'
$ws.Range('C3').Value = 'Ready for handoff'
$ws.Range('D3').Value = '2016-03-25 10:54:50'

# Rebuild hyperlinks: same target addresses per cell position, refreshed display text
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/a840fa4063da3878f03b34b1c495149e0bdc911d/e2e/2e932acd-e47f-4f3c-8372-e61745a5bd03.md', $null, $null, '77232830-5d71-4781-bf09-c9d381f516af.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/a840fa4063da3878f03b34b1c495149e0bdc911d/e2e/77232830-5d71-4781-bf09-c9d381f516af.md', $null, $null, '2e932acd-e47f-4f3c-8372-e61745a5bd03.md') | Out-Null

# ---- Sheet: zh-cn ----
$ws = $wb.Worksheets.Item('zh-cn')

# Update cell values (row 2 <-> row 3 content swap + refreshed status/dates)
$ws.Range('A2').Value = '77232830-5d71-4781-bf09-c9d381f516af.md'
$ws.Range('B2').Value = '.md'
$ws.Range('C2').Value = 'Handed back: in sync with en-US'
$ws.Range('D2').Value = '77232830-5d71-4781-bf09-c9d381f516af.0e8e17f28b30abe6247ce0fbbee84489a351d336.zh-cn.xlf'
$ws.Range('E2').Value = '2016-03-25 10:54:40'
$ws.Range('F2').Value = '77232830-5d71-4781-bf09-c9d381f516af.md'
$ws.Range('G2').Value = '77232830-5d71-4781-bf09-c9d381f516af.0e8e17f28b30abe6247ce0fbbee84489a351d336.zh-cn.xlf'
$ws.Range('H2').Value = '2016-03-25 10:53:43'
$ws.Range('J2').Value = 'Include'
$ws.Range('A3').Value = '2e932acd-e47f-4f3c-8372-e61745a5bd03.md'
$ws.Range('B3').Value = '.md'
$ws.Range('C3').Value = 'Ready for handoff'
$ws.Range('D3').Value = '2e932acd-e47f-4f3c-8372-e61745a5bd03.5ee2385c5afa136df3464c6653ed3e6a3161a3d8.zh-cn.xlf'
$ws.Range('E3').Value = '2016-03-25 10:54:40'
$ws.Range('F3').Value = '2e932acd-e47f-4f3c-8372-e61745a5bd03.md'
$ws.Range('G3').Value = '2e932acd-e47f-4f3c-8372-e61745a5bd03.5ee2385c5afa136df3464c6653ed3e6a3161a3d8.zh-cn.xlf'
$ws.Range('H3').Value = '2016-03-25 10:53:43'
$ws.Range('J3').Value = 'Include'

# Rebuild hyperlinks: same target addresses per cell position, refreshed display text
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/a840fa4063da3878f03b34b1c495149e0bdc911d/e2e/2e932acd-e47f-4f3c-8372-e61745a5bd03.md', $null, $null, '77232830-5d71-4781-bf09-c9d381f516af.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('D2'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/94b2af89cfdad9fae4b802d0f17b1861ba7d8d6c/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/2e932acd-e47f-4f3c-8372-e61745a5bd03.5ee2385c5afa136df3464c6653ed3e6a3161a3d8.zh-cn.xlf', $null, $null, '77232830-5d71-4781-bf09-c9d381f516af.0e8e17f28b30abe6247ce0fbbee84489a351d336.zh-cn.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F2'), 'https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/859af12d6ef2320c741e8736b6c0c53a0f5f8d77/e2e/2e932acd-e47f-4f3c-8372-e61745a5bd03.md', $null, $null, '77232830-5d71-4781-bf09-c9d381f516af.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('G2'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/bc31554d1e2eea4ed1b4726d4a94e5f93d6ec413/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/2e932acd-e47f-4f3c-8372-e61745a5bd03.5ee2385c5afa136df3464c6653ed3e6a3161a3d8.zh-cn.xlf', $null, $null, '77232830-5d71-4781-bf09-c9d381f516af.0e8e17f28b30abe6247ce0fbbee84489a351d336.zh-cn.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/a840fa4063da3878f03b34b1c495149e0bdc911d/e2e/77232830-5d71-4781-bf09-c9d381f516af.md', $null, $null, '2e932acd-e47f-4f3c-8372-e61745a5bd03.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('D3'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/94b2af89cfdad9fae4b802d0f17b1861ba7d8d6c/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/77232830-5d71-4781-bf09-c9d381f516af.0e8e17f28b30abe6247ce0fbbee84489a351d336.zh-cn.xlf', $null, $null, '2e932acd-e47f-4f3c-8372-e61745a5bd03.5ee2385c5afa136df3464c6653ed3e6a3161a3d8.zh-cn.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F3'), 'https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/859af12d6ef2320c741e8736b6c0c53a0f5f8d77/e2e/77232830-5d71-4781-bf09-c9d381f516af.md', $null, $null, '2e932acd-e47f-4f3c-8372-e61745a5bd03.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('G3'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/bc31554d1e2eea4ed1b4726d4a94e5f93d6ec413/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/77232830-5d71-4781-bf09-c9d381f516af.0e8e17f28b30abe6247ce0fbbee84489a351d336.zh-cn.xlf', $null, $null, '2e932acd-e47f-4f3c-8372-e61745a5bd03.5ee2385c5afa136df3464c6653ed3e6a3161a3d8.zh-cn.xlf') | Out-Null

# ---- Sheet: de-de ----
$ws = $wb.Worksheets.Item('de-de')

# Update cell values (row 2 <-> row 3 content swap + refreshed status/dates)
$ws.Range('A2').Value = '77232830-5d71-4781-bf09-c9d381f516af.md'
$ws.Range('B2').Value = '.md'
$ws.Range('C2').Value = 'Handed back: in sync with en-US'
$ws.Range('D2').Value = '77232830-5d71-4781-bf09-c9d381f516af.0e8e17f28b30abe6247ce0fbbee84489a351d336.de-de.xlf'
$ws.Range('E2').Value = '2016-03-25 10:54:50'
$ws.Range('F2').Value = '77232830-5d71-4781-bf09-c9d381f516af.md'
$ws.Range('G2').Value = '77232830-5d71-4781-bf09-c9d381f516af.0e8e17f28b30abe6247ce0fbbee84489a351d336.de-de.xlf'
$ws.Range('H2').Value = '2016-03-25 10:53:58'
$ws.Range('J2').Value = 'Include'
$ws.Range('A3').Value = '2e932acd-e47f-4f3c-8372-e61745a5bd03.md'
$ws.Range('B3').Value = '.md'
$ws.Range('C3').Value = 'Ready for handoff'
$ws.Range('D3').Value = '2e932acd-e47f-4f3c-8372-e61745a5bd03.5ee2385c5afa136df3464c6653ed3e6a3161a3d8.de-de.xlf'
$ws.Range('E3').Value = '2016-03-25 10:54:50'
$ws.Range('F3').Value = '2e932acd-e47f-4f3c-8372-e61745a5bd03.md'
$ws.Range('G3').Value = '2e932acd-e47f-4f3c-8372-e61745a5bd03.5ee2385c5afa136df3464c6653ed3e6a3161a3d8.de-de.xlf'
$ws.Range('H3').Value = '2016-03-25 10:53:58'
$ws.Range('J3').Value = 'Include'

# Rebuild hyperlinks: same target addresses per cell position, refreshed display text
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/a840fa4063da3878f03b34b1c495149e0bdc911d/e2e/2e932acd-e47f-4f3c-8372-e61745a5bd03.md', $null, $null, '77232830-5d71-4781-bf09-c9d381f516af.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('D2'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9a0097503ccc873554da1958355484159f060a44/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/2e932acd-e47f-4f3c-8372-e61745a5bd03.5ee2385c5afa136df3464c6653ed3e6a3161a3d8.de-de.xlf', $null, $null, '77232830-5d71-4781-bf09-c9d381f516af.0e8e17f28b30abe6247ce0fbbee84489a351d336.de-de.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F2'), 'https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/d8491c8e638ce102003153519de200c3718ee99b/e2e/2e932acd-e47f-4f3c-8372-e61745a5bd03.md', $null, $null, '77232830-5d71-4781-bf09-c9d381f516af.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('G2'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/82202019ef099c35c84e427084ef3558fe381345/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/2e932acd-e47f-4f3c-8372-e61745a5bd03.5ee2385c5afa136df3464c6653ed3e6a3161a3d8.de-de.xlf', $null, $null, '77232830-5d71-4781-bf09-c9d381f516af.0e8e17f28b30abe6247ce0fbbee84489a351d336.de-de.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/a840fa4063da3878f03b34b1c495149e0bdc911d/e2e/77232830-5d71-4781-bf09-c9d381f516af.md', $null, $null, '2e932acd-e47f-4f3c-8372-e61745a5bd03.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('D3'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9a0097503ccc873554da1958355484159f060a44/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/77232830-5d71-4781-bf09-c9d381f516af.0e8e17f28b30abe6247ce0fbbee84489a351d336.de-de.xlf', $null, $null, '2e932acd-e47f-4f3c-8372-e61745a5bd03.5ee2385c5afa136df3464c6653ed3e6a3161a3d8.de-de.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F3'), 'https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/d8491c8e638ce102003153519de200c3718ee99b/e2e/77232830-5d71-4781-bf09-c9d381f516af.md', $null, $null, '2e932acd-e47f-4f3c-8372-e61745a5bd03.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('G3'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/82202019ef099c35c84e427084ef3558fe381345/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/77232830-5d71-4781-bf09-c9d381f516af.0e8e17f28b30abe6247ce0fbbee84489a351d336.de-de.xlf', $null, $null, '2e932acd-e47f-4f3c-8372-e61745a5bd03.5ee2385c5afa136df3464c6653ed3e6a3161a3d8.de-de.xlf') | Out-Null

